# Insert a new weekly price record as row 174 ("Haba", Femacal de La Calera),
# shifting the existing rows 174-223 down to 175-224.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 174 (pushes old row 174..223 down to 175..224)
$ws.Rows("174").Insert()

# Fill in the new row 174 with the new record's data
$ws.Range("A174").Value = 3
$ws.Range("B174").Value = "Femacal de La Calera"
$ws.Range("C174").Value = "Coquimbo"
$ws.Range("D174").Value = 44876
$ws.Range("E174").Value = 5
$ws.Range("F174").Value = 100112026
$ws.Range("G174").Value = "Haba"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 75
$ws.Range("K174").Value = 8000
$ws.Range("L174").Value = 8500
$ws.Range("M174").Value = 8267
$ws.Range("N174").Value = "$/saco 25 kilos"
$ws.Range("O174").Value = "Provincia de Limarí"
$ws.Range("P174").Value = 331
$ws.Range("Q174").Value = 25
$ws.Range("R174").Value = "Hortaliza"
